$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-25 06:33:33"
$wsZhCn.Range("G3").Value = "2016-01-25 06:34:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-25 06:33:44"
$wsDeDe.Range("G3").Value = "2016-01-25 06:34:37"
